$wb = $excel.ActiveWorkbook

# Existing sheets
$channel = $wb.Worksheets.Item("Channel")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Add a new "Discover" worksheet after the last existing sheet (Channel)
$discover = $wb.Worksheets.Add($null, $lastSheet)
$discover.Name = "Discover"

# Build the target header style (yellow fill reused from Channel's header,
# general/bottom alignment, no wrap) on a scratch cell first, so the new
# cellXfs entry ends up identical to the one the sheet header cells need.
$channel.Range("A1").Copy()
$discover.Range("Z1").PasteSpecial(-4122)
$discover.Range("Z1").HorizontalAlignment = 1
$discover.Range("Z1").VerticalAlignment = -4107

# Apply that exact style to the real header range, then clear the scratch cell
$discover.Range("Z1").Copy()
$discover.Range("A1:D1").PasteSpecial(-4122)
$discover.Range("Z1").Clear()

# Header row values
$discover.Range("A1").Value = "Count"
$discover.Range("B1").Value = "Channel Name"
$discover.Range("C1").Value = "Channel type"
$discover.Range("D1").Value = "Followers"

# Match the selection captured in the saved file
$discover.Range("B2").Select() | Out-Null
